$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 881, shifting existing rows 881:922 down to 882:923
$ws.Rows.Item(881).Insert()

# Populate the newly inserted row with the new data point
# Force column A to remain plain text (not auto-converted to a date serial)
$ws.Cells.Item(881, 1).NumberFormat = "@"
$ws.Cells.Item(881, 1).Value = "2026/02/25"
$ws.Cells.Item(881, 2).Value = "水"
$ws.Cells.Item(881, 3).Value = 13
$ws.Cells.Item(881, 4).Value = 201
